# Form the consolidated report: fill in the "Absent" (column H) values
# that had not yet been computed/consolidated for a few rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value  = 1
$ws.Range("H4").Value  = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 0
